# Rotate the varying data (Fecha, Volumen, Precio mínimo/máximo/promedio,
# Origen, Precio $/Kg) among rows 2, 3 and 4 of the "Pepino ensalada" sheet:
#   new row2 = old row3
#   new row3 = old row4
#   new row4 = old row2
# Columns A,B,C,E,F,G,H,I,N,Q,R are identical across the three rows, so only
# D, J, K, L, M, O, P need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Capture original values for rows 2-4 before overwriting anything.
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{
        2 = $ws.Range("${col}2").Value2
        3 = $ws.Range("${col}3").Value2
        4 = $ws.Range("${col}4").Value2
    }
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $orig[$col][3]
    $ws.Range("${col}3").Value2 = $orig[$col][4]
    $ws.Range("${col}4").Value2 = $orig[$col][2]
}
